$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.478.20'
$ws.Range("E2").Value = '  -1.33%  '
$ws.Range("D3").Value = '1.849.01'
$ws.Range("E3").Value = '  -0.62%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9995'
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '243.13'
$ws.Range("E5").Value = '  -1.38%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6590'
$ws.Range("E6").Value = '  +3.58%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.000'
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '47.99'
$ws.Range("E8").Value = '  +2.95%  '
$ws.Range("B9").Value = 'Cardano'
$ws.Range("C9").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2992'
$ws.Range("E9").Value = '  -0.24%  '
$ws.Range("B10").Value = 'Dogecoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07487'
$ws.Range("E10").Value = '  +0.11%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '24.35'
$ws.Range("E11").Value = '  -0.57%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07640'
$ws.Range("D13").Value = '1.841.71'
$ws.Range("E13").Value = '  -0.94%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.022'
$ws.Range("E14").Value = '  -0.59%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6846'
$ws.Range("E15").Value = '  -0.67%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '83.75'
$ws.Range("E16").Value = '  -0.72%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000009537'
$ws.Range("E17").Value = '  +1.42%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.146'
$ws.Range("E18").Value = '  +0.77%  '
$ws.Range("D19").Value = '29.518.94'
$ws.Range("E19").Value = '  -1.04%  '
$ws.Range("D20").Value = '2.078.47'
$ws.Range("E20").Value = '  -1.49%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '236.96'
$ws.Range("E21").Value = '  -0.81%  '
$ws.Range("E22").Value = '  -0.64%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.9996'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.671'
$ws.Range("E24").Value = '  +4.29%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.001'
$ws.Range("E25").Value = '  -0.15%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1425'
$ws.Range("E26").Value = '  +0.53%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '156.69'
$ws.Range("E27").Value = '  -1.61%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.491'
$ws.Range("E28").Value = '  -1.04%  '
$ws.Range("E29").Value = '  -0.95%  '
$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.490'
$ws.Range("E30").Value = '  -0.90%  '
$ws.Range("B31").Value = 'Hedera'
$ws.Range("C31").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.06027'
$ws.Range("E31").Value = '  -0.62%  '
$ws.Range("E32").Value = '  -1.90%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.137'
$ws.Range("E33").Value = '  -0.15%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.075'
$ws.Range("E34").Value = '  -1.59%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.180'
$ws.Range("E35").Value = '  +1.23%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.854'
$ws.Range("E36").Value = '  -0.69%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.7229'
$ws.Range("E37").Value = '  -0.78%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.595'
$ws.Range("E38").Value = '  -0.76%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.802'
$ws.Range("E39").Value = '  -1.89%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01781'
$ws.Range("E40").Value = '  -0.82%  '
$ws.Range("D41").Value = '1.199.20'
$ws.Range("E41").Value = '  -1.87%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.239'
$ws.Range("E42").Value = '  -1.00%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9091'
$ws.Range("E43").Value = '  -2.12%  '
$ws.Range("E44").Value = '  -0.26%  '
$ws.Range("D45").Value = '2.011.39'
$ws.Range("E45").Value = '  -0.48%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '101.95'
$ws.Range("E46").Value = '  -0.40%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '66.07'
$ws.Range("E47").Value = '  -0.35%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.447'
$ws.Range("E48").Value = '  +10.81%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.4060'
$ws.Range("E49").Value = '  -0.70%  '
$ws.Range("E50").Value = '  -2.90%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '9.046'
$ws.Range("E51").Value = '  -2.47%  '
